# "error solve ifrs list" - fix bad/erroneous IFRS financial figures in the
# company_list sheet: correct the (previously mis-scaled) numeric data for
# fiscal years 2014-2018 (rows 2-6), and clear the erroneous forecast rows
# for 2019/12(E)-2021/12(E) (rows 7-9), leaving only the period labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 957
$ws.Range("E2").Value = -97
$ws.Range("F2").Value = -113
$ws.Range("G2").Value = -125
$ws.Range("H2").Value = -146
$ws.Range("I2").Value = -146
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 968
$ws.Range("L2").Value = 602
$ws.Range("M2").Value = 366
$ws.Range("N2").Value = 364
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 305
$ws.Range("Q2").Value = -63
$ws.Range("R2").Value = 25
$ws.Range("S2").Value = -37
$ws.Range("T2").Value = 29
$ws.Range("U2").Value = -92
$ws.Range("V2").Value = 462
$ws.Range("W2").Value = -10.14
$ws.Range("X2").Value = -15.28
$ws.Range("Y2").Value = -33.34
$ws.Range("Z2").Value = -13.53
$ws.Range("AA2").Value = 164.69
$ws.Range("AB2").Value = 42.06
$ws.Range("AC2").Value = -2392
$ws.Range("AD2").Value = -1.29
$ws.Range("AE2").Value = 6547
$ws.Range("AF2").Value = 0.47
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 6106308

# Row 3
$ws.Range("D3").Value = 890
$ws.Range("E3").Value = -58
$ws.Range("F3").Value = -58
$ws.Range("G3").Value = -124
$ws.Range("H3").Value = -162
$ws.Range("I3").Value = -162
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 747
$ws.Range("L3").Value = 481
$ws.Range("M3").Value = 266
$ws.Range("N3").Value = 265
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 364
$ws.Range("Q3").Value = 13
$ws.Range("R3").Value = 57
$ws.Range("S3").Value = -73
$ws.Range("T3").Value = 11
$ws.Range("U3").Value = 2
$ws.Range("V3").Value = 341
$ws.Range("W3").Value = -6.56
$ws.Range("X3").Value = -18.18
$ws.Range("Y3").Value = -51.41
$ws.Range("Z3").Value = -18.86
$ws.Range("AA3").Value = 180.42
$ws.Range("AB3").Value = -7.91
$ws.Range("AC3").Value = -2616
$ws.Range("AD3").Value = -2.35
$ws.Range("AE3").Value = 3942
$ws.Range("AF3").Value = 1.56
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 7270242

# Row 4
$ws.Range("D4").Value = 576
$ws.Range("E4").Value = -77
$ws.Range("F4").Value = -57
$ws.Range("G4").Value = -182
$ws.Range("H4").Value = -161
$ws.Range("I4").Value = -163
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 792
$ws.Range("L4").Value = 449
$ws.Range("M4").Value = 344
$ws.Range("N4").Value = 299
$ws.Range("O4").Value = 45
$ws.Range("P4").Value = 419
$ws.Range("Q4").Value = -22
$ws.Range("R4").Value = -186
$ws.Range("S4").Value = 215
$ws.Range("T4").Value = 13
$ws.Range("U4").Value = -35
$ws.Range("V4").Value = 328
$ws.Range("W4").Value = -13.46
$ws.Range("X4").Value = -28.01
$ws.Range("Y4").Value = -57.71
$ws.Range("Z4").Value = -20.95
$ws.Range("AA4").Value = 130.64
$ws.Range("AB4").Value = -11.99
$ws.Range("AC4").Value = -2040
$ws.Range("AD4").Value = -2.33
$ws.Range("AE4").Value = 3818
$ws.Range("AF4").Value = 1.24
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 8376352

# Row 5
$ws.Range("D5").Value = 346
$ws.Range("E5").Value = -56
$ws.Range("F5").Value = -56
$ws.Range("G5").Value = -180
$ws.Range("H5").Value = -57
$ws.Range("I5").Value = -59
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 640
$ws.Range("L5").Value = 403
$ws.Range("M5").Value = 237
$ws.Range("N5").Value = 236
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 419
$ws.Range("Q5").Value = -54
$ws.Range("R5").Value = 40
$ws.Range("S5").Value = 13
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = -58
$ws.Range("V5").Value = 279
$ws.Range("W5").Value = -16.24
$ws.Range("X5").Value = -16.57
$ws.Range("Y5").Value = -22.04
$ws.Range("Z5").Value = -8.01
$ws.Range("AA5").Value = 170.02
$ws.Range("AB5").Value = -26.17
$ws.Range("AC5").Value = -704
$ws.Range("AD5").Value = -4.45
$ws.Range("AE5").Value = 3012
$ws.Range("AF5").Value = 1.04
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 8376352

# Row 6
$ws.Range("D6").Value = 360
$ws.Range("E6").Value = -60
$ws.Range("F6").Value = -60
$ws.Range("G6").Value = 47
$ws.Range("H6").Value = 15
$ws.Range("I6").Value = 15
$ws.Range("K6").Value = 725
$ws.Range("L6").Value = 482
$ws.Range("M6").Value = 243
$ws.Range("N6").Value = 242
$ws.Range("P6").Value = 419
$ws.Range("Q6").Value = -37
$ws.Range("R6").Value = 163
$ws.Range("S6").Value = -137
$ws.Range("T6").Value = 4
$ws.Range("U6").Value = -41
$ws.Range("V6").Value = 294
$ws.Range("W6").Value = -16.76
$ws.Range("X6").Value = 4.09
$ws.Range("Y6").Value = 6.16
$ws.Range("Z6").Value = 2.16
$ws.Range("AA6").Value = 197.98
$ws.Range("AB6").Value = -25.27
$ws.Range("AC6").Value = 176
$ws.Range("AD6").Value = 15.75
$ws.Range("AE6").Value = 3094
$ws.Range("AF6").Value = 0.9
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 8376352

# Rows 7-9: clear all data columns except A, B, C
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
